# Graphs side by side with their responses
# Fill in the previously-empty response cells for subject/row 9 ("data" sheet)
# so the plot_* (graph) responses sit alongside the other survey answers,
# and mark the "total_crt" (AL) flag for that row as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
$ws.Activate()

# --- Row 9 : newly answered columns -----------------------------------
$ws.Range("B9").Value  = "yes"
$ws.Range("C9").Value  = "No"
$ws.Range("D9").Value  = "No"
$ws.Range("E9").Value  = "Algo familiarizado"
$ws.Range("F9").Value  = 10
$ws.Range("G9").Value  = 100
$ws.Range("H9").Value  = 47

$ws.Range("K9").Value  = "Poco analítico/a"
$ws.Range("L9").Value  = "No"
$ws.Range("M9").Value  = "Si"

$ws.Range("O9").Value  = "Raramente"
$ws.Range("P9").Value  = 0
$ws.Range("Q9").Value  = "Muy importantes"
$ws.Range("R9").Value  = "Algo interesado/a"
$ws.Range("S9").Value  = "Si"

$ws.Range("V9").Value  = "q"
$ws.Range("W9").Value  = "q"
$ws.Range("X9").Value  = "q"
$ws.Range("Y9").Value  = "p"
$ws.Range("Z9").Value  = "p"
$ws.Range("AA9").Value = "q"
$ws.Range("AB9").Value = "p"
$ws.Range("AC9").Value = "q"
$ws.Range("AD9").Value = "q"
$ws.Range("AE9").Value = "p"
$ws.Range("AF9").Value = "p"
$ws.Range("AG9").Value = "p"
$ws.Range("AH9").Value = "q"
$ws.Range("AI9").Value = "p"
$ws.Range("AJ9").Value = "Menor que la media del mercado"

# total_crt flag flips from 0 to 1 now that the row is fully answered
$ws.Range("AL9").Value = 1

# --- View state: selection moves to AK22 -------------------------------
$ws.Range("AK22").Select()
